$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2, column A: product name changes from "arduino mega" to "arduino uno"
$ws.Range("A2").Value = "arduino uno"

# Row 2, columns B & C: cost and number change (formula in D2 recalculates automatically)
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 4

# Update the active selection to match the final state
$ws.Range("C5").Select()
